$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 124: mark hidden, and flip E124/F124 from 0 to 1 ---
$ws.Range("E124").Value = 1
$ws.Range("F124").Value = 1
$ws.Rows.Item(124).Hidden = $true

# --- Row 126: mark hidden (values unchanged) ---
$ws.Rows.Item(126).Hidden = $true

# --- Insert a new hidden, empty "spacer" row at 127 (copy formatting from row above) ---
$ws.Rows.Item(127).Insert(-4121, 0)
$ws.Range("A127").Clear()
$ws.Rows.Item(127).Hidden = $true

# --- Insert a new data row at 128 (copy formatting from row above) ---
$ws.Rows.Item(128).Insert(-4121, 0)
$ws.Range("A128").WrapText = $true

$ws.Range("A128").Value = "ايجاد، حذف، نمايش و ويرايش بخش آيتم‌هاي تمريني الگوي تمرين"
$ws.Range("B128").Value = $ws.Range("B124").Value()
$ws.Range("C128").Value = 1
$ws.Range("D128").Value = 1
$ws.Range("E128").Value = 0
$ws.Range("F128").Value = 0

# --- Grow the table to cover the two new rows ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F128"))

# --- Update the selected cell to match the saved view ---
$null = $ws.Range("F106").Select()

Write-Host "edit complete"
